# Home Screen test cases are added
#
# Refresh the login credentials used by the "MA_AccountEdit1" test case
# (new admin account + password), record a "pass" expected-result value
# for that scenario, mirror the "Pass" outcome back onto the "Test Cases"
# sheet, and leave the workbook's view state where a user would be after
# making these edits (back on the "Test Cases" tab).

$wb = $excel.ActiveWorkbook
$wsTest = $wb.Worksheets.Item("Test Cases")
$wsAcct = $wb.Worksheets.Item("MA_AccountEdit1")

# --- MA_AccountEdit1 sheet: new login credentials -------------------------
# Username (A2) keeps its existing mailto hyperlink/style; only the
# displayed text is refreshed.
$wsAcct.Range("A2").Value = "admin@nfhslearn.com"

# Password (B2) is a brand-new hyperlinked value, so add the hyperlink and
# then match A2's look (instead of Excel's default blue/underlined link
# style).
$wsAcct.Range("B2").Value = "nfhslearn@6186"
$wsAcct.Hyperlinks.Add($wsAcct.Range("B2"), "mailto:nfhslearn@6186")
$wsAcct.Range("B2").Style = $wsAcct.Range("A2").Style

# New "pass" expected-result value in column D (plain, unstyled cell)
$wsAcct.Range("D2").Value = "pass"
$wsAcct.Range("D2").Style = "Normal"

# Widen the password column so the longer value fits
$wsAcct.Columns.Item(2).ColumnWidth = 30

# --- Test Cases sheet: mirror the "Pass" result ----------------------------
$wsTest.Range("E2").Value = "Pass"

# --- Selection / active sheet bookkeeping ----------------------------------
# Remember MA_AccountEdit1's last selection (C5) while it's still active ...
$wsAcct.Range("C5").Select()
# ... then switch back to Test Cases as the active/selected tab, at B4.
$wsTest.Activate()
$wsTest.Range("B4").Select()
